# Apply the latest cryptos-list price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. thousand-dot grouping like
# "70.143.81", or trailing zeros like "13.30") that must stay TEXT, not be
# reinterpreted as a number. A leading apostrophe forces Excel to store it as text,
# exactly like typing it in by hand.

$ws.Range("D2").Value = "'70.143.81"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "'3.622.79"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'605.02"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'196.84"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'54.04"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.0000306"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "'9.57"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'4.195.56"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "'13.30"
$ws.Range("E15").Value = "  +5.57%  "
$ws.Range("D16").Value = "'591.43"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "'19.25"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'70.281.34"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'3.605.81"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "'0.995"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'17.76"
$ws.Range("E22").Value = "  -3.22%  "
$ws.Range("D23").Value = "'5.16"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'102.48"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").Value = "'3.06"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "'10.79"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'9.64"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").Value = "'33.88"
$ws.Range("D30").Value = "'4.53"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'12.36"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'63.25"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +12.01%  "
$ws.Range("D36").Value = "'3.958.95"
$ws.Range("E36").Value = "  +5.98%  "
$ws.Range("D37").Value = "'3.17"
$ws.Range("E37").Value = "  +5.58%  "
$ws.Range("D38").Value = "'528.02"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'37.39"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "'2.87"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "'3.35"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "'8.63"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  +4.90%  "
$ws.Range("D51").Value = "'1.31"
$ws.Range("E51").Value = "  +3.85%  "
